$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last data row (row 22) so the sheet dimension becomes A1:H21
$ws.Rows.Item(22).Delete()

$newValues = @(
  @(-0.9255759716033934, 1.071177214384079, -0.401322513818741, 0.2278527319431305, -0.8228355050086975, 0.3602577745914459),
  @(-1.280259013175964, 0.5940434336662286, -1.134497284889223, 0.0354301854968071, -0.427452951669693, 0.4167627990245819),
  @(-1.571854650974273, 0.08382582664489754, -2.215982675552368, 0.1412625908851623, -0.3335324227809906, -0.0261144898831844),
  @(-1.946246147155765, 0.2504835128784156, -2.480054140090942, 0.152563601732254, -0.527634859085083, -0.0335975885391235),
  @(-5.073675036430361, -1.826430678367616, -2.739734649658203, 0.2588541507720947, -0.7525860071182251, 0.0536034256219863),
  @(-3.199407458305356, -1.501654863357543, -3.262988328933716, 0.447000652551651, 0.88361656665802, -0.780380368232727),
  @(8.78546768426898, -8.89101099967958, 47.06318786740315, 2.5057692527771, -4.185495853424072, -1.749212622642517),
  @(3.72559547424312, -6.610608488321275, 35.79121446609481, -3.3768630027771, -6.853449821472168, 1.152244567871094),
  @(-8.249110460281372, -1.285275377333166, -3.368247583508488, -3.452763080596924, -1.31610918045044, -1.093448758125305),
  @(-7.083846211433405, -2.286428689956659, 0.4840306788682932, 0.9560040831565856, 1.219134330749511, 3.047301292419434),
  @(-1.840081989765169, 2.111219294369222, 0.04846413433551781, -0.0424551330506801, 0.1565342247486114, 0.2597704529762268),
  @(3.811614036560053, 0.9647219777107299, 0.292347490787505, 0.1169806942343711, -0.6998988389968872, -1.625512361526489),
  @(1.02576005458831, 3.592965073883541, -0.1755489408969892, 1.398270487785339, -1.411556839942932, -0.2585487067699432),
  @(0.2673106193542476, 1.979196399450299, 0.2698189914226538, 2.174527645111084, 1.24815046787262, 0.5600107312202454),
  @(0.6833226680755626, -0.3385796248912853, -1.02221550792456, -0.1776090711355209, 0.5535966157913208, -0.1540907621383667),
  @(0.4996232986450184, -0.3156740069389323, -0.2882215976715059, -0.0699440687894821, 0.0386372283101081, -0.0519235469400882),
  @(-0.2515395879745499, 0.7405525892972948, 1.177142247557639, -0.0510072484612464, 0.0113010071218013, -0.0332921557128429),
  @(-0.430771350860593, 0.8799855709075934, 0.8326354324817646, 0.0404698215425014, 0.3057380616664886, 0.0070249503478407),
  @(0.2168139219284046, 1.084875062108039, 0.5181054919958111, -0.0203112699091434, 0.0247400421649217, -0.1175915598869323),
  @(-0.08790016174316459, 1.008570432662963, 0.4307903051376349, -0.0080939643085002, -0.0186313893646001, -0.0047342055477201)
)

$r = 2
foreach ($rowVals in $newValues) {
    $c = 3
    foreach ($val in $rowVals) {
        $ws.Cells.Item($r, $c).Value2 = $val
        $c = $c + 1
    }
    $r = $r + 1
}

"Done: dimension=" + $ws.UsedRange.Address()